# Update the "Förändrad" (Changed) date column (C) for rows 2-171
# from Excel serial date 45202 (2023-10-03) to 45203 (2023-10-04).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C171").Value = 45203
